$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename: column C now reports "emissions" instead of the old
# "Target Max Electricity kWh per anum" label (column D keeps "Star Rating"). ---
$ws.Range("C1").Value = "emissions"

# --- Update the column C (emissions) figures for all 4 "Star" rating rows ---
# New values = old Target-Max-Electricity figures scaled by the emissions
# factor (~1.09) used to convert kWh to a scope emissions number.
$ws.Range("C2").Value = 79915.1
$ws.Range("C3").Value = 93871.3
$ws.Range("C4").Value = 107827.5
$ws.Range("C5").Value = 159830.2
$ws.Range("C6").Value = 187742.6
$ws.Range("C7").Value = 215654.9
$ws.Range("C8").Value = 239745.3
$ws.Range("C9").Value = 281613.9
$ws.Range("C10").Value = 323482.4
$ws.Range("C11").Value = 319660.4
$ws.Range("C12").Value = 375485.2
$ws.Range("C13").Value = 431309.9
$ws.Range("C14").Value = 399575.5
$ws.Range("C15").Value = 469356.4
$ws.Range("C16").Value = 539137.4

# --- Re-format column C: drop the bespoke Verdana/grey font, switch the
# number format from integer-with-thousands (#,##0) to 2-decimal
# (#,##0.00) since the new emissions figures carry fractional kWh. ---
$cRange = $ws.Range("C2:C16")
$cRange.ClearFormats()
$cRange.NumberFormat = "#,##0.00"

# --- Selection moved from H9 to F2 ---
$ws.Range("F2").Select()
